$wb = $excel.ActiveWorkbook

# --- Sheet: Schedule (rows 2-5, cols A-F) ---
$ws1 = $wb.Worksheets.Item("Schedule")
$schedule = @(
  @(46081.3125, 46081.8125, 12.0, 45.36, 231.1974405, 5.096945337301587),
  @(46081.83333333334, 46082.0, 4.0, 15.12, 400.3730835, 26.47970128968254),
  @(46082.04166666666, 46082.20833333334, 4.0, 15.12, 289.66275, 19.15758928571429),
  @(46082.3125, 46082.64583333334, 8.0, 30.24, 87.933339, 2.907848511904762)
)
$r = 2
foreach ($row in $schedule) {
  $c = 1
  foreach ($val in $row) {
    $ws1.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}

# --- Sheet: Detailed (rows 2-97, cols A-E) ---
$ws2 = $wb.Worksheets.Item("Detailed")
$detailed = @(
  @(46081.0, 65.0, "historical", 46081.0, "OFF"),
  @(46081.02083333334, 66.73901, "historical", 46081.0, "OFF"),
  @(46081.04166666666, 57.06, "historical", 46081.0, "OFF"),
  @(46081.0625, 57.06, "historical", 46081.0, "OFF"),
  @(46081.08333333334, 56.98, "historical", 46081.0, "OFF"),
  @(46081.10416666666, 56.14298, "historical", 46081.0, "OFF"),
  @(46081.125, 55.3107, "historical", 46081.0, "OFF"),
  @(46081.14583333334, 54.99126, "historical", 46081.0, "OFF"),
  @(46081.16666666666, 55.00664, "historical", 46081.0, "OFF"),
  @(46081.1875, 55.52054, "historical", 46081.0, "OFF"),
  @(46081.20833333334, 56.1745, "historical", 46081.0, "OFF"),
  @(46081.22916666666, 65.0, "historical", 46081.0, "OFF"),
  @(46081.25, 57.36, "historical", 46081.0, "OFF"),
  @(46081.27083333334, 57.36, "historical", 46081.0, "OFF"),
  @(46081.29166666666, 56.98, "historical", 46081.0, "OFF"),
  @(46081.3125, 35.88, "historical", 46081.0, "ON"),
  @(46081.33333333334, 1.16971, "historical", 46081.0, "ON"),
  @(46081.35416666666, 1.15893, "historical", 46081.0, "ON"),
  @(46081.375, 1.16138, "historical", 46081.0, "ON"),
  @(46081.39583333334, 1.13838, "historical", 46081.0, "ON"),
  @(46081.41666666666, 1.07591, "historical", 46081.0, "ON"),
  @(46081.4375, 0.7, "historical", 46081.0, "ON"),
  @(46081.45833333334, 0.51, "historical", 46081.0, "ON"),
  @(46081.47916666666, 0.7, "historical", 46081.0, "ON"),
  @(46081.5, 0.50995, "historical", 46081.0, "ON"),
  @(46081.52083333334, 0.51, "historical", 46081.0, "ON"),
  @(46081.54166666666, 0.7, "historical", 46081.0, "ON"),
  @(46081.5625, 0.7, "historical", 46081.0, "ON"),
  @(46081.58333333334, 0.7, "historical", 46081.0, "ON"),
  @(46081.60416666666, 1.11062, "historical", 46081.0, "ON"),
  @(46081.625, 0.51, "historical", 46081.0, "ON"),
  @(46081.64583333334, 0.51, "historical", 46081.0, "ON"),
  @(46081.66666666666, -3.76, "historical", 46081.0, "ON"),
  @(46081.6875, -4.14539, "historical", 46081.0, "ON"),
  @(46081.70833333334, 27.67198, "historical", 46081.0, "ON"),
  @(46081.72916666666, 35.88, "historical", 46081.0, "ON"),
  @(46081.75, 38.7, "historical", 46081.0, "ON"),
  @(46081.77083333334, 45.77354, "forecast", 46081.0, "ON"),
  @(46081.79166666666, 48.26057, "forecast", 46081.0, "ON"),
  @(46081.8125, 53.08559, "forecast", 46081.0, "OFF"),
  @(46080.83333333334, 74.22528, "historical", 46080.0, "OFF"),
  @(46080.85416666666, 78.0, "historical", 46080.0, "OFF"),
  @(46080.875, 71.4, "historical", 46080.0, "OFF"),
  @(46080.89583333334, 65.0, "historical", 46080.0, "OFF"),
  @(46080.91666666666, 64.89, "historical", 46080.0, "OFF"),
  @(46080.9375, 65.0, "historical", 46080.0, "OFF"),
  @(46080.95833333334, 64.1059, "historical", 46080.0, "OFF"),
  @(46080.97916666666, 65.0, "forecast", 46081.0, "ON"),
  @(46082.0, 37.89, "forecast", 46082.0, "OFF"),
  @(46082.02083333334, 52.93539, "forecast", 46082.0, "OFF"),
  @(46082.04166666666, 37.89, "forecast", 46082.0, "ON"),
  @(46082.0625, 37.89, "forecast", 46082.0, "ON"),
  @(46082.08333333334, 37.89, "forecast", 46082.0, "ON"),
  @(46082.10416666666, 35.88, "forecast", 46082.0, "ON"),
  @(46082.125, 35.88, "forecast", 46082.0, "ON"),
  @(46082.14583333334, 35.88, "forecast", 46082.0, "ON"),
  @(46082.16666666666, 37.89, "forecast", 46082.0, "ON"),
  @(46082.1875, 37.89, "forecast", 46082.0, "ON"),
  @(46082.20833333334, 37.89, "forecast", 46082.0, "OFF"),
  @(46082.22916666666, 37.89, "forecast", 46082.0, "OFF"),
  @(46082.25, 37.89, "forecast", 46082.0, "OFF"),
  @(46082.27083333334, 37.89, "forecast", 46082.0, "OFF"),
  @(46082.29166666666, 37.89, "forecast", 46082.0, "OFF"),
  @(46082.3125, 0.7, "forecast", 46082.0, "ON"),
  @(46082.33333333334, 0.51, "forecast", 46082.0, "ON"),
  @(46082.35416666666, 0.7, "forecast", 46082.0, "ON"),
  @(46082.375, 0.7, "forecast", 46082.0, "ON"),
  @(46082.39583333334, 0.51, "forecast", 46082.0, "ON"),
  @(46082.41666666666, 0.01133, "forecast", 46082.0, "ON"),
  @(46082.4375, 0.7, "forecast", 46082.0, "ON"),
  @(46082.45833333334, 0.7, "forecast", 46082.0, "ON"),
  @(46082.47916666666, 2.01653, "forecast", 46082.0, "ON"),
  @(46082.5, 1.92022, "forecast", 46082.0, "ON"),
  @(46082.52083333334, 16.37624, "forecast", 46082.0, "ON"),
  @(46082.54166666666, 7.97027, "forecast", 46082.0, "ON"),
  @(46082.5625, 5.59048, "forecast", 46082.0, "ON"),
  @(46082.58333333334, 8.88867, "forecast", 46082.0, "ON"),
  @(46082.60416666666, 17.73128, "forecast", 46082.0, "ON"),
  @(46082.625, 25.16302, "forecast", 46082.0, "ON"),
  @(46082.64583333334, 37.89, "forecast", 46082.0, "OFF"),
  @(46082.66666666666, 37.89, "forecast", 46082.0, "OFF"),
  @(46082.6875, 39.6168, "forecast", 46082.0, "OFF"),
  @(46082.70833333334, 48.98861, "forecast", 46082.0, "OFF"),
  @(46082.72916666666, 51.63038, "forecast", 46082.0, "OFF"),
  @(46082.75, 57.06009, "forecast", 46082.0, "OFF"),
  @(46082.77083333334, 57.31, "forecast", 46082.0, "OFF"),
  @(46082.79166666666, 57.31, "forecast", 46082.0, "OFF"),
  @(46082.8125, 57.31, "forecast", 46082.0, "OFF"),
  @(46082.83333333334, 57.31, "forecast", 46082.0, "OFF"),
  @(46082.85416666666, 57.0601, "forecast", 46082.0, "OFF"),
  @(46082.875, 57.06, "forecast", 46082.0, "OFF"),
  @(46082.89583333334, 57.06, "forecast", 46082.0, "OFF"),
  @(46082.91666666666, 51.40578, "forecast", 46082.0, "OFF"),
  @(46082.9375, 56.98, "forecast", 46082.0, "OFF"),
  @(46082.95833333334, 51.93775, "forecast", 46082.0, "OFF"),
  @(46082.97916666666, 51.94373, "forecast", 46082.0, "OFF")
)
$r = 2
foreach ($row in $detailed) {
  $c = 1
  foreach ($val in $row) {
    $ws2.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}
